$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: the automatic "_GoBack" bookmark moves away from
#    here (it will be re-created at the point of the last text edit,
#    inside the "Demandas..." bullet, below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) "Demandas ... via Trello ..." bullet: text stays the same, but it
#    was retyped around "via" (leaving the cursor -- and so the new
#    "_GoBack" bookmark -- right after "vi").
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" serão criadas via Trello ", $true, $false, $false, $false, $false, $true, 1, $false, " serão criadas via Trello ", 2)

$r2 = $d.Content
$r2.Find.Execute("serão criadas vi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPoint = $d.Range($r2.End, $r2.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# ------------------------------------------------------------------
# 3) "Duvidas sempre ..." bullet: add the missing accent -> "Dúvidas".
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Duvidas", $true, $false, $false, $false, $false, $true, 1, $false, "Dúvidas", 2)

# ------------------------------------------------------------------
# 4) "Cliente deve Acompanhar ..." bullet: fix the capitalisation of
#    "Acompanhar", "FeedBacks", "Escopo" and "Desenvolvimento".
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Cliente deve Acompanhar e gerar FeedBacks durante Escopo e Desenvolvimento.", $true, $false, $false, $false, $false, $true, 1, $false, "Cliente deve acompanhar e gerar feedbacks durante escopo e desenvolvimento.", 2)
